# Auto-generated edit script: update Response/Processed/Category/Explanation
# columns for rows that previously failed (timeout) or had outdated text,
# per commit 'feat: new data exp25'.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = ' Category: CAT1
Explanation: The incident description involves multiple failed authentication attempts from an external IP address on various user accounts (test5, rundeck, jimmy, and root). This is indicative of account compromise attempts, such as credential phishing or brute force attacks, which fall under the CAT1 category.'
$ws.Range("B2").Value = '{''Category'': ''CAT1'', ''Explanation'': ''The incident description involves multiple failed authentication attempts from an external IP address on various user accounts (test5, rundeck, jimmy, and root). This is indicative of account compromise attempts, such as credential phishing or brute force attacks, which fall under the CAT1 category.''}'
$ws.Range("C2").Value = 'CAT1'
$ws.Range("D2").Value = 'The incident description involves multiple failed authentication attempts from an external IP address on various user accounts (test5, rundeck, jimmy, and root). This is indicative of account compromise attempts, such as credential phishing or brute force attacks, which fall under the CAT1 category.'

$ws.Range("A3").Value = ' Category: CAT1
Explanation: The incident description involves multiple failed authentication attempts from a single IP address, which suggests an attempt at account compromise through brute force attacks. This aligns with the definition of CAT1: Account Compromise.'
$ws.Range("B3").Value = '{''Category'': ''CAT1'', ''Explanation'': ''The incident description involves multiple failed authentication attempts from a single IP address, which suggests an attempt at account compromise through brute force attacks. This aligns with the definition of CAT1: Account Compromise.''}'
$ws.Range("C3").Value = 'CAT1'
$ws.Range("D3").Value = 'The incident description involves multiple failed authentication attempts from a single IP address, which suggests an attempt at account compromise through brute force attacks. This aligns with the definition of CAT1: Account Compromise.'

$ws.Range("A4").Value = ' Category: CAT7
Explanation: The incident description involves an email reporting abuse from a specific IP address, which was detected by an automated process. This is a clear case of social engineering as the sender is attempting to deceive the recipient into investigating and taking action based on the information provided. This falls under the category of phishing, a common form of social engineering attack.'
$ws.Range("B4").Value = '{''Category'': ''CAT7'', ''Explanation'': ''The incident description involves an email reporting abuse from a specific IP address, which was detected by an automated process. This is a clear case of social engineering as the sender is attempting to deceive the recipient into investigating and taking action based on the information provided. This falls under the category of phishing, a common form of social engineering attack.''}'
$ws.Range("C4").Value = 'CAT7'
$ws.Range("D4").Value = 'The incident description involves an email reporting abuse from a specific IP address, which was detected by an automated process. This is a clear case of social engineering as the sender is attempting to deceive the recipient into investigating and taking action based on the information provided. This falls under the category of phishing, a common form of social engineering attack.'

$ws.Range("A5").Value = ' Category: CAT1
Explanation: The incident description involves multiple failed SSH login attempts from an external IP address that does not belong to the organization. This could be a brute force attack aimed at compromising user accounts, which falls under the Account Compromise (CAT1) category.'
$ws.Range("B5").Value = '{''Category'': ''CAT1'', ''Explanation'': ''The incident description involves multiple failed SSH login attempts from an external IP address that does not belong to the organization. This could be a brute force attack aimed at compromising user accounts, which falls under the Account Compromise (CAT1) category.''}'
$ws.Range("C5").Value = 'CAT1'
$ws.Range("D5").Value = 'The incident description involves multiple failed SSH login attempts from an external IP address that does not belong to the organization. This could be a brute force attack aimed at compromising user accounts, which falls under the Account Compromise (CAT1) category.'

$ws.Range("A6").Value = ' Category: CAT2
Explanation: The incident description involves a malware infection by SystemBC, which is a known type of malicious code. This matches the definition of category CAT2: Malware.'
$ws.Range("B6").Value = '{''Category'': ''CAT2'', ''Explanation'': ''The incident description involves a malware infection by SystemBC, which is a known type of malicious code. This matches the definition of category CAT2: Malware.''}'
$ws.Range("C6").Value = 'CAT2'
$ws.Range("D6").Value = 'The incident description involves a malware infection by SystemBC, which is a known type of malicious code. This matches the definition of category CAT2: Malware.'

$ws.Range("A7").Value = ' Category: CAT2
Explanation: The incident description mentions a malware infection by SystemBC, which is a type of malicious code. This matches the definition of Malware (CAT2) in the NIST categories.'
$ws.Range("B7").Value = '{''Category'': ''CAT2'', ''Explanation'': ''The incident description mentions a malware infection by SystemBC, which is a type of malicious code. This matches the definition of Malware (CAT2) in the NIST categories.''}'
$ws.Range("C7").Value = 'CAT2'
$ws.Range("D7").Value = 'The incident description mentions a malware infection by SystemBC, which is a type of malicious code. This matches the definition of Malware (CAT2) in the NIST categories.'

$ws.Range("A8").Value = ' Category: CAT11
Explanation: The incident involves a third-party issue as it originates from an organization (ORGANIZATION_dc2075266f) reporting abuse from an IP address (IP_ADDRESS_984299dc8f), which is on the network of another organization (pfsense-svp.furg.br). This suggests a security incident involving suppliers or service providers, as per CAT11.'
$ws.Range("B8").Value = '{''Category'': ''CAT11'', ''Explanation'': ''The incident involves a third-party issue as it originates from an organization (ORGANIZATION_dc2075266f) reporting abuse from an IP address (IP_ADDRESS_984299dc8f), which is on the network of another organization (pfsense-svp.furg.br). This suggests a security incident involving suppliers or service providers, as per CAT11.''}'
$ws.Range("C8").Value = 'CAT11'
$ws.Range("D8").Value = 'The incident involves a third-party issue as it originates from an organization (ORGANIZATION_dc2075266f) reporting abuse from an IP address (IP_ADDRESS_984299dc8f), which is on the network of another organization (pfsense-svp.furg.br). This suggests a security incident involving suppliers or service providers, as per CAT11.'

$ws.Range("A9").Value = ' Category: CAT11
Explanation: The incident involves a third-party issue as it concerns an abuse report from another organization ([ORGANIZATION_dc2075266f]) regarding the IP address [IP_ADDRESS_984299dc8f], which is on the network of the recipient organization. This indicates that the security incident originated from a supplier or service provider, falling under CAT11.'
$ws.Range("B9").Value = '{''Category'': ''CAT11'', ''Explanation'': ''The incident involves a third-party issue as it concerns an abuse report from another organization ([ORGANIZATION_dc2075266f]) regarding the IP address [IP_ADDRESS_984299dc8f], which is on the network of the recipient organization. This indicates that the security incident originated from a supplier or service provider, falling under CAT11.''}'
$ws.Range("C9").Value = 'CAT11'
$ws.Range("D9").Value = 'The incident involves a third-party issue as it concerns an abuse report from another organization ([ORGANIZATION_dc2075266f]) regarding the IP address [IP_ADDRESS_984299dc8f], which is on the network of the recipient organization. This indicates that the security incident originated from a supplier or service provider, falling under CAT11.'

$ws.Range("A10").Value = ' Category: CAT1
Explanation: The incident description mentions unauthorized access to a network (possibly compromised machine) and its use for malicious activities. This aligns with the definition of Account Compromise (CAT1).'
$ws.Range("B10").Value = '{''Category'': ''CAT1'', ''Explanation'': ''The incident description mentions unauthorized access to a network (possibly compromised machine) and its use for malicious activities. This aligns with the definition of Account Compromise (CAT1).''}'
$ws.Range("C10").Value = 'CAT1'
$ws.Range("D10").Value = 'The incident description mentions unauthorized access to a network (possibly compromised machine) and its use for malicious activities. This aligns with the definition of Account Compromise (CAT1).'

$ws.Range("A11").Value = ' Category: CAT3
Explanation: The incident description involves a DDoS (Denial of Service Attack) on one of the organization''s customers using an IP address under their control. This aligns with examples provided for CAT3, such as volumetric DoS or DDoS attacks.'
$ws.Range("B11").Value = '{''Category'': ''CAT3'', ''Explanation'': "The incident description involves a DDoS (Denial of Service Attack) on one of the organization''s customers using an IP address under their control. This aligns with examples provided for CAT3, such as volumetric DoS or DDoS attacks."}'
$ws.Range("C11").Value = 'CAT3'
$ws.Range("D11").Value = 'The incident description involves a DDoS (Denial of Service Attack) on one of the organization''s customers using an IP address under their control. This aligns with examples provided for CAT3, such as volumetric DoS or DDoS attacks.'

$ws.Range("A12").Value = ' Category: CAT3
Explanation: The incident description involves a Denial of Service Attack as it mentions a DDoS attack on a specific IP address, which makes systems unavailable.'
$ws.Range("B12").Value = '{''Category'': ''CAT3'', ''Explanation'': ''The incident description involves a Denial of Service Attack as it mentions a DDoS attack on a specific IP address, which makes systems unavailable.''}'
$ws.Range("C12").Value = 'CAT3'
$ws.Range("D12").Value = 'The incident description involves a Denial of Service Attack as it mentions a DDoS attack on a specific IP address, which makes systems unavailable.'

$ws.Range("A13").Value = ' Category: CAT3
Explanation: The incident description involves a Distributed Denial of Service (DDoS) attack on a specific IP address, which falls under the category of Denial of Service Attack (CAT3). The attack was part of a coordinated botnet and caused significant packet loss for the affected clients.'
$ws.Range("B13").Value = '{''Category'': ''CAT3'', ''Explanation'': ''The incident description involves a Distributed Denial of Service (DDoS) attack on a specific IP address, which falls under the category of Denial of Service Attack (CAT3). The attack was part of a coordinated botnet and caused significant packet loss for the affected clients.''}'
$ws.Range("C13").Value = 'CAT3'
$ws.Range("D13").Value = 'The incident description involves a Distributed Denial of Service (DDoS) attack on a specific IP address, which falls under the category of Denial of Service Attack (CAT3). The attack was part of a coordinated botnet and caused significant packet loss for the affected clients.'

$ws.Range("A14").Value = ' Category: CAT5
Explanation: The incident description details a vulnerability in the Zimbra Collaboration Suite (CVE-2022-37042) that allows an attacker to upload files, execute arbitrary code, and gain unauthorized access to user accounts. This is a clear example of Vulnerability Exploitation (CAT5).'
$ws.Range("B14").Value = '{''Category'': ''CAT5'', ''Explanation'': ''The incident description details a vulnerability in the Zimbra Collaboration Suite (CVE-2022-37042) that allows an attacker to upload files, execute arbitrary code, and gain unauthorized access to user accounts. This is a clear example of Vulnerability Exploitation (CAT5).''}'
$ws.Range("C14").Value = 'CAT5'
$ws.Range("D14").Value = 'The incident description details a vulnerability in the Zimbra Collaboration Suite (CVE-2022-37042) that allows an attacker to upload files, execute arbitrary code, and gain unauthorized access to user accounts. This is a clear example of Vulnerability Exploitation (CAT5).'

$ws.Range("A15").Value = ' Category: CAT3
Explanation: The incident description details open BGP (179/tcp) services exposed to the internet on multiple IP addresses, which can potentially lead to Denial of Service attacks on this port and affect the availability of the BGP service for the affected system autonomous. This aligns with CAT3: Denial of Service Attack as it involves making systems unavailable due to a flood of requests or traffic.'
$ws.Range("B15").Value = '{''Category'': ''CAT3'', ''Explanation'': ''The incident description details open BGP (179/tcp) services exposed to the internet on multiple IP addresses, which can potentially lead to Denial of Service attacks on this port and affect the availability of the BGP service for the affected system autonomous. This aligns with CAT3: Denial of Service Attack as it involves making systems unavailable due to a flood of requests or traffic.''}'
$ws.Range("C15").Value = 'CAT3'
$ws.Range("D15").Value = 'The incident description details open BGP (179/tcp) services exposed to the internet on multiple IP addresses, which can potentially lead to Denial of Service attacks on this port and affect the availability of the BGP service for the affected system autonomous. This aligns with CAT3: Denial of Service Attack as it involves making systems unavailable due to a flood of requests or traffic.'

$ws.Range("A16").Value = ' Category: CAT12
Explanation: The incident description involves an IP address performing unwanted activities towards a server and being blocked by the system. This is indicative of an intrusion attempt, as the IP address is attempting to connect to various ports on the server, which could be part of a port scan or other malicious activity. However, since it''s not confirmed whether the attack was successful, it falls under CAT12: Intrusion Attempt.'
$ws.Range("B16").Value = '{''Category'': ''CAT12'', ''Explanation'': "The incident description involves an IP address performing unwanted activities towards a server and being blocked by the system. This is indicative of an intrusion attempt, as the IP address is attempting to connect to various ports on the server, which could be part of a port scan or other malicious activity. However, since it''s not confirmed whether the attack was successful, it falls under CAT12: Intrusion Attempt."}'
$ws.Range("C16").Value = 'CAT12'
$ws.Range("D16").Value = 'The incident description involves an IP address performing unwanted activities towards a server and being blocked by the system. This is indicative of an intrusion attempt, as the IP address is attempting to connect to various ports on the server, which could be part of a port scan or other malicious activity. However, since it''s not confirmed whether the attack was successful, it falls under CAT12: Intrusion Attempt.'

$ws.Range("A17").Value = ' Category: CAT3
Explanation: The incident description details a service (SSDP/UPnP) that can be exploited for DDoS attacks due to its exposure to the internet. This service can be used to amplify responses, causing a denial of service attack on other organizations and consuming more bandwidth.'
$ws.Range("B17").Value = '{''Category'': ''CAT3'', ''Explanation'': ''The incident description details a service (SSDP/UPnP) that can be exploited for DDoS attacks due to its exposure to the internet. This service can be used to amplify responses, causing a denial of service attack on other organizations and consuming more bandwidth.''}'
$ws.Range("C17").Value = 'CAT3'
$ws.Range("D17").Value = 'The incident description details a service (SSDP/UPnP) that can be exploited for DDoS attacks due to its exposure to the internet. This service can be used to amplify responses, causing a denial of service attack on other organizations and consuming more bandwidth.'

$ws.Range("A18").Value = ' Category: CAT3
Explanation: The incident description involves servers under the responsibility of the recipient being open to monlist and readvar queries, which can be abused for distributed denial-of-service (DDoS) attacks. This aligns with the definition of a Denial of Service Attack (CAT3). Examples include volumetric DoS or DDoS (UDP flood, SYN flood, HTTP/HTTPS flood), attacks on APIs or websites, Mirai botnet. The search terms in the description also support this classification: "Denial of service", "distributed denial-of-service", "attaques distribuídos de negação de serviço".'
$ws.Range("B18").Value = '{''Category'': ''CAT3'', ''Explanation'': ''The incident description involves servers under the responsibility of the recipient being open to monlist and readvar queries, which can be abused for distributed denial-of-service (DDoS) attacks. This aligns with the definition of a Denial of Service Attack (CAT3). Examples include volumetric DoS or DDoS (UDP flood, SYN flood, HTTP/HTTPS flood), attacks on APIs or websites, Mirai botnet. The search terms in the description also support this classification: "Denial of service", "distributed denial-of-service", "attaques distribuídos de negação de serviço".''}'
$ws.Range("C18").Value = 'CAT3'
$ws.Range("D18").Value = 'The incident description involves servers under the responsibility of the recipient being open to monlist and readvar queries, which can be abused for distributed denial-of-service (DDoS) attacks. This aligns with the definition of a Denial of Service Attack (CAT3). Examples include volumetric DoS or DDoS (UDP flood, SYN flood, HTTP/HTTPS flood), attacks on APIs or websites, Mirai botnet. The search terms in the description also support this classification: "Denial of service", "distributed denial-of-service", "attaques distribuídos de negação de serviço".'

$ws.Range("A19").Value = ' Category: CAT7
Explanation: The incident description involves a phishing scam, which falls under the Social Engineering (CAT7) category. The email contains deceptive content and attempts to trick the recipient into revealing sensitive information or clicking on malicious links.'
$ws.Range("B19").Value = '{''Category'': ''CAT7'', ''Explanation'': ''The incident description involves a phishing scam, which falls under the Social Engineering (CAT7) category. The email contains deceptive content and attempts to trick the recipient into revealing sensitive information or clicking on malicious links.''}'
$ws.Range("C19").Value = 'CAT7'
$ws.Range("D19").Value = 'The incident description involves a phishing scam, which falls under the Social Engineering (CAT7) category. The email contains deceptive content and attempts to trick the recipient into revealing sensitive information or clicking on malicious links.'

$ws.Range("A20").Value = ' Category: CAT2
Explanation: The incident description mentions an email server on the network that is sending fraudulent emails attempting to extort users. This strongly suggests malware infection, specifically a type of spam or phishing campaign, aligning with the Malware (CAT2) category.'
$ws.Range("B20").Value = '{''Category'': ''CAT2'', ''Explanation'': ''The incident description mentions an email server on the network that is sending fraudulent emails attempting to extort users. This strongly suggests malware infection, specifically a type of spam or phishing campaign, aligning with the Malware (CAT2) category.''}'
$ws.Range("C20").Value = 'CAT2'
$ws.Range("D20").Value = 'The incident description mentions an email server on the network that is sending fraudulent emails attempting to extort users. This strongly suggests malware infection, specifically a type of spam or phishing campaign, aligning with the Malware (CAT2) category.'

$ws.Range("A22").Value = ' Category: CAT9
Explanation: The incident description states that the content of a website has been altered without authorization, which falls under the category of unauthorized modification.'
$ws.Range("B22").Value = '{''Category'': ''CAT9'', ''Explanation'': ''The incident description states that the content of a website has been altered without authorization, which falls under the category of unauthorized modification.''}'
$ws.Range("C22").Value = 'CAT9'
$ws.Range("D22").Value = 'The incident description states that the content of a website has been altered without authorization, which falls under the category of unauthorized modification.'

$ws.Range("A23").Value = ' Category: CAT9
Explanation: The incident description indicates an unauthorized modification of a website''s content, which falls under the category of Unauthorized Modification (CAT9).'
$ws.Range("B23").Value = '{''Category'': ''CAT9'', ''Explanation'': "The incident description indicates an unauthorized modification of a website''s content, which falls under the category of Unauthorized Modification (CAT9)."}'
$ws.Range("C23").Value = 'CAT9'
$ws.Range("D23").Value = 'The incident description indicates an unauthorized modification of a website''s content, which falls under the category of Unauthorized Modification (CAT9).'

$ws.Range("A25").Value = ' Category: CAT9
Explanation: The incident description shows that an unauthorized party has altered the content of a website, which falls under the category of Unauthorized Modification (CAT9). The hacker group ''theMx0nday'' has taken credit for the defacement and left their signature on the site.'
$ws.Range("B25").Value = '{''Category'': ''CAT9'', ''Explanation'': "The incident description shows that an unauthorized party has altered the content of a website, which falls under the category of Unauthorized Modification (CAT9). The hacker group ''theMx0nday'' has taken credit for the defacement and left their signature on the site."}'
$ws.Range("C25").Value = 'CAT9'
$ws.Range("D25").Value = 'The incident description shows that an unauthorized party has altered the content of a website, which falls under the category of Unauthorized Modification (CAT9). The hacker group ''theMx0nday'' has taken credit for the defacement and left their signature on the site.'

Write-Host "Applied updates to rows: 2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,25"
